$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" positioned between "2021-Q4" and "总计"
# ---------------------------------------------------------------------------
$sheetQ4 = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($null, $sheetQ4)
$newSheet.Name = "2022-Q1"

# Copy the row formatting (borders / bold / alignment) used on the other
# quarter sheets (header row style + first-column "index" style) onto the
# same region of the new sheet.
$sheetQ4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$sheetQ4.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# --- header row -------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- data rows ----------------------------------------------------------
$data = @(
    @(0, "001128", "宝盈新兴产业灵活配置混合",   "19.94", "86.76", "3.40", "0.6780", 10),
    @(1, "519678", "银河消费驱动混合",           "1.06",  "75.49", "7.78", "0.0825", 4),
    @(2, "011431", "泰达宏利消费服务混合A",       "1.61",  "81.15", "4.01", "0.0646", 1),
    @(3, "011073", "鹏华安润混合A",              "3.07",  "29.79", "1.01", "0.0310", 5),
    @(4, "011432", "泰达宏利消费服务混合C",       "0.15",  "81.15", "4.01", "0.0060", 1),
    @(5, "011074", "鹏华安润混合C",              "0.25",  "29.79", "1.01", "0.0025", 5)
)

# Force text storage for the fund-code / percentage-looking columns so that
# values such as "001128" keep their leading zero instead of becoming the
# number 1128.
$dataBlock = $newSheet.Range("B2:G7")
$dataBlock.NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Re-apply the plain formatting on top of the values: this clears the
# temporary "text" number-format created above (so the saved file has no
# stray style index on these cells) while leaving the stored value type
# (text vs number) untouched.
$sheetQ4.Range("C2").Copy()
$dataBlock.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add a row for 2022-Q1 above the
#    existing rows.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.86

# renumber the running index in column A for the rows that shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

Write-Host "2022-Q1 sheet inserted and 总计 sheet updated"
